$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

# New model row data (mirrors the layout of the existing rows).
$ws.Cells.Item($row, 1).Value = "llama-3-8b-bnb-4bit-synthetic_text_to_sql-lora-3epochs-Q5_K_M:latest"
$ws.Cells.Item($row, 2).Value = "llama3:70b"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 200
$ws.Cells.Item($row, 5).Value = 2615.66

# Columns F:AH (except AI) and AJ:AP have no data for this run - write them
# as blank/empty text cells (matching the rest of the sheet's convention of
# empty string placeholders rather than leaving the cells completely unset).
$blank1 = $ws.Range("F" + $row + ":AH" + $row)
$blank1.Value = "'"
$blank1.ClearFormats()

$ws.Cells.Item($row, 35).Value = 2056.28

$blank2 = $ws.Range("AJ" + $row + ":AP" + $row)
$blank2.Value = "'"
$blank2.ClearFormats()

$ws.Cells.Item($row, 43).Value = 559.38
$ws.Cells.Item($row, 44).Value = 70
$ws.Cells.Item($row, 45).Value = "logs\llama_3_8b_bnb_4bit_synthetic_text_to_sql_lora_3epochs_Q5_K_M_latest_llama3_70b_1_200_test_bootstrap_match_1.txt"
$ws.Cells.Item($row, 46).Value = 559.38
$ws.Cells.Item($row, 47).Value = 35
$ws.Cells.Item($row, 48).Value = "logs\llama_3_8b_bnb_4bit_synthetic_text_to_sql_lora_3epochs_Q5_K_M_latest_llama3_70b_1_200_test_bootstrap_correct_1.txt"
$ws.Cells.Item($row, 49).Value = 58.33333333333334
$ws.Cells.Item($row, 50).Value = 4
$ws.Cells.Item($row, 51).Value = 8
